# Update countries & provincias Spain
# Applies the 26-May-2020 data refresh (11:05 -> 11:35) to the "Pais" sheet:
#  - updates the "last updated" timestamp banner
#  - refreshes case counters for India, Belgica, Rumania(ctx), Austria,
#    Marruecos/Malasia (which swap table order), Albania(ctx), and the
#    Trinidad y Tobago / Siria / Malaui trio (which also swap order)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header / timestamp banner -------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 26 de Mayo de 2020 a las 11:35"

# --- Row 13 (India) --------------------------------------------------------
$ws.Range("B13").Value = 146208
$ws.Range("C13").Value = 1258
$ws.Range("D13").Value = 61052
$ws.Range("E13").Value = 80969
$ws.Range("F13").Value = 0
$ws.Range("G13").Value = 15
$ws.Range("H13").Value = 4187

# --- Row 22 (Belgica) -------------------------------------------------------
$ws.Range("B22").Value = 57455
$ws.Range("C22").Value = 113
$ws.Range("D22").Value = 15320
$ws.Range("E22").Value = 32801
$ws.Range("F22").Value = 0
$ws.Range("G22").Value = 22
$ws.Range("H22").Value = 9334

# --- Row 40 -------------------------------------------------------------
$ws.Range("E40").Value = 5443
$ws.Range("G40").Value = 5
$ws.Range("H40").Value = 1210

# --- Row 44 (Austria) -----------------------------------------------------
$ws.Range("B44").Value = 16557
$ws.Range("C44").Value = 18
$ws.Range("D44").Value = 15182
$ws.Range("E44").Value = 732
$ws.Range("F44").Value = 0
$ws.Range("G44").Value = 2
$ws.Range("H44").Value = 643

# --- Rows 60/61: Marruecos & Malasia swap places ---------------------------
# Row 60 used to be Marruecos, it now becomes Malasia with fresh totals;
# row 61 used to be Malasia, it now becomes Marruecos carrying over the
# figures that used to sit in row 60.
$ws.Range("A60").Value = "Malasia"
$ws.Range("B60").Value = 7604
$ws.Range("C60").Value = 187
$ws.Range("D60").Value = 6041
$ws.Range("E60").Value = 1448
$ws.Range("F60").Value = 0
$ws.Range("G60").Value = 0
$ws.Range("H60").Value = 115

$ws.Range("A61").Value = "Marruecos"
$ws.Range("B61").Value = 7532
$ws.Range("C61").Value = 0
$ws.Range("D61").Value = 4774
$ws.Range("E61").Value = 2558
$ws.Range("F61").Value = 0
$ws.Range("G61").Value = 0
$ws.Range("H61").Value = 200

# --- Row 113 (Mali) ----------------------------------------------------
$ws.Range("B113").Value = 1029
$ws.Range("C113").Value = 25
$ws.Range("D113").Value = 803
$ws.Range("E113").Value = 193
$ws.Range("F113").Value = 0
$ws.Range("G113").Value = 1
$ws.Range("H113").Value = 33

# --- Rows 168/169: Trinidad y Tobago & Siria swap places --------------------
# Row 168 used to be Trinidad yTobago, it now becomes Siria with fresh
# totals; row 169 used to be Siria, it now becomes Trinidad yTobago
# carrying over the figures that used to sit in row 168.
$ws.Range("A168").Value = "Siria"
$ws.Range("B168").Value = 121
$ws.Range("C168").Value = 15
$ws.Range("D168").Value = 41
$ws.Range("E168").Value = 76
$ws.Range("F168").Value = 0
$ws.Range("G168").Value = 0
$ws.Range("H168").Value = 4

$ws.Range("A169").Value = "Trinidad yTobago"
$ws.Range("B169").Value = 116
$ws.Range("C169").Value = 0
$ws.Range("D169").Value = 108
$ws.Range("E169").Value = 0
$ws.Range("F169").Value = 0
$ws.Range("G169").Value = 0
$ws.Range("H169").Value = 8

# --- Row 170 (Malaui) ----------------------------------------------------
$ws.Range("D170").Value = 37
$ws.Range("E170").Value = 60
